# DALA-5788: Adjusted Excel Files for Alias Export
# Insert a new "Alias Export" column right after "Field Name" (column D),
# duplicating the Field Name values, and shifting all subsequent columns
# (old E:L -> new F:M) to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remember the width of column D (Field Name) so the newly inserted column
# can be given a matching width, mirroring how Excel copies the left
# neighbour's column formatting on a normal column insert.
$fieldNameWidth = $ws.Columns("D").ColumnWidth

# Insert a new column before column E. This shifts existing columns E:L to
# F:M and, like a normal Excel column insert, the new column inherits the
# formatting of the column to its left (D).
$ws.Columns("E:E").Insert()
$ws.Columns("E").ColumnWidth = $fieldNameWidth

# Header for the newly inserted column.
$ws.Range("E1").Value = "Alias Export"

# Duplicate "Field Name" (column D) values into the new "Alias Export"
# column (E) for every data row.
$lastRow = $ws.Cells($ws.Rows.Count, "D").End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    if ($null -ne $dVal) {
        $ws.Cells.Item($r, 5).Value = $dVal
    }
}

# Keep the active selection on E6, matching the author's saved view state.
$ws.Range("E6").Select()
